$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.188.54"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "1.905.80"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.80"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5250"
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3783"
$ws.Range("E8").Value = "  +3.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07272"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8963"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07686"
$ws.Range("E12").Value = "  +2.51%  "
$ws.Range("D13").Value = "1.912.44"
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.98"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.274"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008621"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.48"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "27.253.43"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.069"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").Value = "2.149.33"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.434"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.310"
$ws.Range("E25").Value = "  +10.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.81"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.735"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.78"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.800"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09225"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8108"
$ws.Range("E33").Value = "  +8.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05054"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.239"
$ws.Range("E35").Value = "  +7.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.986"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.309"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.587"
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5681"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01987"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.982"
$ws.Range("E42").Value = "  +5.14%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.21"
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.618"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1514"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4829"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.23"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.621"
$ws.Range("E49").Value = "  +4.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.56"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.72"
$ws.Range("E51").Value = "  +0.95%  "
